$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dbUrl = "https://www.amazon.com/Willwin-Female-Accessories-Connector-Adaptor/dp/B074QGFDM8/ref=sr_1_7_sspa?crid=4TJRCFDRAYB1&dchild=1&keywords=db9%2Bconnector&qid=1596222297&sprefix=DB9%2Caps%2C171&sr=8-7-spons&spLa=ZW5jcnlwdGVkUXVhbGlmaWVyPUEzS084M1pCTUZXNkRFJmVuY3J5cHRlZElkPUExMDI5Njk5Mjg3TjgxV0FaRUlQVSZlbmNyeXB0ZWRBZElkPUEwODgxOTg4MzlHWVoxUkNOUVhLUCZ3aWRnZXROYW1lPXNwX210ZiZhY3Rpb249Y2xpY2tSZWRpcmVjdCZkb05vdExvZ0NsaWNrPXRydWU&th=1"
$dbUrlShort = $dbUrl.Substring(0, 255)
$mcUrl = "https://www.mcmaster.com/95117A499/"

# --- Quantity updates on existing rows ---
$ws.Range("B12").Value2 = 2
$ws.Range("B21").Value2 = 12
$ws.Range("B22").Value2 = 12

# --- Column width tweaks (auto-fit side effect of new, wider content) ---
$ws.Columns.Item(1).ColumnWidth = 13.43
$ws.Columns.Item(3).ColumnWidth = 21.43

# --- New rows 24-27 ---
$ws.Range("A24").Value2 = "Interface Box"
$ws.Range("B24").Value2 = 1
$ws.Range("C24").Value2 = "DB9 Connector Female"
$ws.Range("D24").Value2 = "Willwin"
$ws.Range("E24").Value2 = "Willwin DB9 Female Connector"
$ws.Range("F24").Value2 = $dbUrl
$ws.Range("F24").Style = "Hyperlink"

$ws.Range("A25").Value2 = "Sensor Mounts"
$ws.Range("B25").Value2 = 1
$ws.Range("C25").Value2 = "DB9 Connector Male"
$ws.Range("D25").Value2 = "Willwin"
$ws.Range("E25").Value2 = "Willwin DB9 Male Connector"
$ws.Range("F25").Value2 = $dbUrl
$ws.Range("F25").Style = "Hyperlink"

$ws.Range("A26").Value2 = "Interface Box"
$ws.Range("B26").Value2 = 2
$ws.Range("C26").Value2 = "M3 Broaching Nuts"
$ws.Range("F26").Value2 = $mcUrl
$ws.Range("F26").Style = "Hyperlink"

$ws.Range("A27").Value2 = "Interface Box"
$ws.Range("B27").Value2 = 2
$ws.Range("C27").Value2 = "M3 x 6mm Screws"

# --- Hyperlinks ---
# F10 / F4 previously had "empty" hyperlink entries (no r:id); give them a
# display-only (local) hyperlink, same as before, but now carrying the
# truncated display text.
$ws.Hyperlinks.Add($ws.Range("F10"), "", "", "", $dbUrlShort) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "", "", "", $dbUrlShort) | Out-Null

# New rows: F24/F25 are also display-only (local) hyperlinks; F26 is a real
# external hyperlink relationship.
$ws.Hyperlinks.Add($ws.Range("F24"), "", "", "", $dbUrlShort) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F25"), "", "", "", $dbUrlShort) | Out-Null
$ws.Hyperlinks.Add($ws.Range("F26"), $mcUrl) | Out-Null

# Re-apply the Hyperlink cell style (Hyperlinks.Add mints a fresh xf record
# internally; make sure the visible cells end up on the canonical style).
$ws.Range("F10").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F24").Style = "Hyperlink"
$ws.Range("F25").Style = "Hyperlink"
$ws.Range("F26").Style = "Hyperlink"

# --- Selection cursor position ---
$ws.Range("D33").Select() | Out-Null

Write-Output "done"
